$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D12").Value = -6.056099999999997
$ws.Range("D32").Value = -7.305499999999995
$ws.Range("D36").Value = -7.323300000000002
$ws.Range("D38").Value = -7.214299999999998
$ws.Range("D46").Value = -8.126599999999994
$ws.Range("D54").Value = -7.911800000000002
$ws.Range("D55").Value = -7.356799999999995
$ws.Range("D67").Value = -7.373499999999996
$ws.Range("D69").Value = -7.453699999999999
$ws.Range("D72").Value = -7.363399999999997
$ws.Range("D91").Value = -7.960799999999997
$ws.Range("D99").Value = -8.023699999999998
$ws.Range("D104").Value = -7.617299999999996
